# Cryptos price/volume refresh (GitHub Actions scheduled update).
# Mirrors the upstream diff: per-row Price (D) / Volume(1h) (E) refresh,
# plus one ranking swap (rows 32/33: NEARProtocol <-> EthereumClassic).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a cell as literal TEXT even when its content parses as a
# number (e.g. "596.14"), matching the inlineStr/text cells already on the
# sheet. A leading apostrophe forces Excel to store it as text; re-applying
# the "Normal" style afterwards clears the resulting quote-prefix flag so
# the cell keeps its original (unstyled) look.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "67.407.66"
$ws.Range("E2").Value = "  -3.29%  "
$ws.Range("D3").Value = "3.694.13"
$ws.Range("E3").Value = "  -3.32%  "
$ws.Range("E4").Value = "  -0.16%  "
Set-TextValue $ws.Range("D5") "596.14"
$ws.Range("E5").Value = "  -2.31%  "
Set-TextValue $ws.Range("D6") "165.22"
$ws.Range("E6").Value = "  -5.13%  "
$ws.Range("D7").Value = "3.691.26"
$ws.Range("E7").Value = "  -3.38%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.32%  "
Set-TextValue $ws.Range("D10") "0.159"
$ws.Range("E10").Value = "  -3.58%  "
$ws.Range("E11").Value = "  -3.62%  "
$ws.Range("E12").Value = "  -3.17%  "
Set-TextValue $ws.Range("D13") "37.52"
$ws.Range("E13").Value = "  -5.83%  "
$ws.Range("E14").Value = "  -4.83%  "
$ws.Range("D15").Value = "4.309.79"
$ws.Range("E15").Value = "  -3.29%  "
$ws.Range("D16").Value = "3.694.22"
$ws.Range("E16").Value = "  -3.39%  "
$ws.Range("D17").Value = "67.444.43"
$ws.Range("E17").Value = "  -3.34%  "
Set-TextValue $ws.Range("D18") "17.55"
$ws.Range("E18").Value = "  +6.10%  "
$ws.Range("E19").Value = "  -3.49%  "
$ws.Range("E20").Value = "  -3.32%  "
Set-TextValue $ws.Range("D21") "491.38"
$ws.Range("E21").Value = "  -2.45%  "
Set-TextValue $ws.Range("D22") "9.16"
$ws.Range("E22").Value = "  -3.15%  "
Set-TextValue $ws.Range("D23") "0.725"
$ws.Range("E23").Value = "  -1.03%  "
Set-TextValue $ws.Range("D24") "85.63"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  -5.70%  "
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("E27").Value = "  -3.17%  "
Set-TextValue $ws.Range("D28") "10.05"
$ws.Range("E28").Value = "  -3.32%  "
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("E31").Value = "  -6.22%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D32") "31.54"
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D33") "7.63"
$ws.Range("E33").Value = "  -3.99%  "
$ws.Range("D34").Value = "3.831.33"
$ws.Range("E34").Value = "  -3.38%  "
Set-TextValue $ws.Range("D35") "0.107"
$ws.Range("E35").Value = "  -4.78%  "
$ws.Range("D36").Value = "3.631.61"
$ws.Range("E36").Value = "  -3.29%  "
Set-TextValue $ws.Range("D37") "0.999"
$ws.Range("E37").Value = "  -0.18%  "
Set-TextValue $ws.Range("D38") "0.992"
$ws.Range("E38").Value = "  -4.47%  "
$ws.Range("E39").Value = "  -4.87%  "
$ws.Range("E40").Value = "  -6.80%  "
$ws.Range("E41").Value = "  -3.80%  "
Set-TextValue $ws.Range("D42") "435.23"
$ws.Range("E42").Value = "  -10.65%  "
Set-TextValue $ws.Range("D43") "48.61"
$ws.Range("E43").Value = "  -2.18%  "
$ws.Range("E44").Value = "  -5.64%  "
$ws.Range("E45").Value = "  -6.54%  "
Set-TextValue $ws.Range("D46") "8.37"
$ws.Range("E46").Value = "  -1.25%  "
Set-TextValue $ws.Range("D48") "40.48"
$ws.Range("E48").Value = "  -6.23%  "
Set-TextValue $ws.Range("D49") "142.16"
$ws.Range("E49").Value = "  +1.62%  "
$ws.Range("D50").Value = "2.752.38"
$ws.Range("E50").Value = "  -5.64%  "
Set-TextValue $ws.Range("D51") "0.0346"
$ws.Range("E51").Value = "  -3.53%  "
